$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Id changes, Ost/Nord change to old row 11's coords, comment (AC10) cleared
$ws.Range("A10").Value = 111528980
$ws.Range("Q10").Value = 467799.8074815667
$ws.Range("R10").Value = 6875539.119922069
$ws.Range("AC10").Value = ""

# Row 11: Id changes to new value, Lokalnamn + coords change to old row 12's values
$ws.Range("A11").Value = 111527876
$ws.Range("P11").Value = "Fläcksberget V, Hjd"
$ws.Range("Q11").Value = 467615.2905344999
$ws.Range("R11").Value = 6875426.740629551

# Row 12: Id, Lokalnamn + coords change to old row 10's values, comment (AC12) added
$ws.Range("A12").Value = 111528300
$ws.Range("P12").Value = "Fläcksberget, Hjd"
$ws.Range("Q12").Value = 467795.2212022893
$ws.Range("R12").Value = 6875452.272210476
$ws.Range("AC12").Value = "Tre blommande."
